$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Preconditions (E), Method Inputs (F), and Expected Result (G)
# columns for each test case row (rows 7-22), completing the test plan
# after initializing and running the BankAccount unit tests.

$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "account_number=123 client_number=456 account_balance=100.0"
$ws.Range("G7").Value = "Object created"
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "account_number=123 client_number=456 account_balance=`"invalid`""
$ws.Range("G8").Value = 0
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "account_number=`"Wendy`" client_number=456 account_balance=100.0"
$ws.Range("G9").Value = "ValueError"
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "account_number=123 client_number=`"Ways`" account_balance=100.0"
$ws.Range("G10").Value = "ValueError"
$ws.Range("E11").Value = "BankAccount(123,456,100.0)"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "account._BankAccount__account_number=123"
$ws.Range("E12").Value = "BankAccount(123,456,100.0)"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "account._BankAccount__client_number=456"
$ws.Range("E13").Value = "BankAccount(123,456,100.0)"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "account._BankAccount__balance=100.0"
$ws.Range("E14").Value = "BankAccount(123,456,100.0)"
$ws.Range("F14").Value = "update_balance(50.0)"
$ws.Range("G14").Value = "account_BankAccount__balance=150.0"
$ws.Range("E15").Value = "BankAccount(123,456,100.0)"
$ws.Range("F15").Value = "update_balance(-40.0)"
$ws.Range("G15").Value = "account_BankAccount__balance=60.0"
$ws.Range("E16").Value = "BankAccount(123,456,100.0)"
$ws.Range("F16").Value = "update_balance(`"invalid`")"
$ws.Range("G16").Value = "account_BankAccount__balance=100.0"
$ws.Range("E17").Value = "BankAccount(123,456,100.0)"
$ws.Range("F17").Value = "deposit(25.0)"
$ws.Range("G17").Value = "account_BankAccount__balance=125.0"
$ws.Range("E18").Value = "BankAccount(123,456,100.0)"
$ws.Range("F18").Value = "deposit(-20)"
$ws.Range("G18").Value = "ValueError"
$ws.Range("E19").Value = "BankAccount(123,456,100.0)"
$ws.Range("F19").Value = "withdraw(50.0)"
$ws.Range("G19").Value = "account_BankAccount__balance=50.0"
$ws.Range("E20").Value = "BankAccount(123,456,100.0)"
$ws.Range("F20").Value = "withdraw(-10)"
$ws.Range("G20").Value = "ValueError"
$ws.Range("E21").Value = "BankAccount(123,456,100.0)"
$ws.Range("F21").Value = "withdraw(200.0)"
$ws.Range("G21").Value = "ValueError"
$ws.Range("E22").Value = "BankAccount(123,456,100.0)"
$ws.Range("F22").Value = "str(aacount)"
$ws.Range("G22").Value = "`"Client 456, Account 123, Balance: 100.0`""

# Leave the selection where the author ended up after entering the data.
$ws.Range("D19").Select()
